$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated forecast values for column C (y_1) and column E (y_1_forecast)
$values = @{
    2  = @{ C = 7.288845540350142;  E = 10.06916370210014 }
    3  = @{ C = -20.40984652067478; E = -33.45158075171641 }
    4  = @{ C = 5.331710924091837;  E = 10.84949786623359 }
    5  = @{ C = 11.04982736891558;  E = 9.589921161142879 }
    6  = @{ C = 4.748210439985256;  E = 7.086193663490992 }
    7  = @{ C = -2.313034291448757; E = -7.020874871669158 }
    8  = @{ C = 3.711391384148;     E = 0.6610955960690834 }
    9  = @{ C = 3.860244074450181;  E = 3.254220449867029 }
    10 = @{ C = 2.370939381494686;  E = 3.967543131851214 }
    11 = @{ C = 4.421855465610292;  E = 5.260364862099642 }
    12 = @{ C = 3.320585727896552;  E = -4.098213472638578 }
    13 = @{ C = 1.782333336406405;  E = 4.060401000000002 }
    14 = @{ C = -4.477718018907028; E = -11.83522404790002 }
    15 = @{ C = 6.317691071509768;  E = 4.613033063261129 }
    16 = @{ C = 3.79744344971964;   E = 3.703837953294542 }
    17 = @{ C = 0.4998689793225486; E = -0.08221002454066317 }
    18 = @{ C = 0.1967053802870877; E = 6.289039804796182 }
    19 = @{ C = -2.056549539789942; E = -0.3884660724497446 }
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row].C
    $ws.Cells.Item($row, 5).Value = $values[$row].E
}
